# Updated some trades and adjusted closing prices
# Refresh the sector-breakdown metrics (rows 2-15) across all sectors (cols B-M)
# with the latest computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per metric row. Columns, in order, are:
# B=Staples, C=Discretionary, D=Energy, E=REITs, F=Financials, G=Healthcare,
# H=Industrials, I=Utilities, J=Macro, K=Technology, L=Media, M=Fixed Income
$newValues = @{
    2  = @(-1.309,-0.676,-13.933,-1.344,-16.013,0.031,-1.005,-0.916,-15.28,-2.927,-5.713,2.088)
    3  = @(0.842,10.053,-16.505,-2.173,-10.709,9.396000000000001,3.25,15.921,-10.514,3.877,4.371,6.853)
    4  = @(-5.837,9.763999999999999,-31.548,-4.019,-14.548,5.155,-32.299,-7.059,-18.107,24.456,3.06,1.22)
    5  = @(-0.983,15.347,-30.345,4.888,-16.005,17.743,-29.583,-3.525,-10.242,29.167,9.747999999999999,13.436)
    6  = @(10.415,101.988,-28.545,71.16500000000001,47.183,32.073,-22.409,26.731,4.022,133.252,77.97799999999999,7.319)
    7  = @(2.021,15.259,-6.564,11.467,8.119999999999999,5.78,-4.996,4.902,0.8,18.658,12.35,1.437)
    8  = @(0.865,0.991,1.088,0.973,1.035,1.014,1.127,0.883,0.899,1.069,0.961,0.612)
    9  = @(-4.989,7.324,-15.206,3.664,-0.133,-2.325,-13.926,-2.243,-6.457,10.15,4.635,-3.716)
    10 = @(0.07099999999999999,0.669,-0.218,0.435,0.315,0.224,-0.207,0.163,0.007,0.623,0.554,0.063)
    11 = @(0.016,0.147,-0.067,0.111,0.07199999999999999,0.05,-0.05,0.048,0.001,0.168,0.122,0.013)
    12 = @(-42.584,-30.202,-64.238,-43.964,-39.97,-36.067,-50.672,-42.943,-40.828,-40.259,-33.848,-29.365)
    13 = @(19.132,21.804,33.155,24.851,23.634,22.85,27.356,25.95,19.723,28.872,21.087,12.17)
    14 = @(0.463,0.662,0.475,0.415,0.707,0.569,0.479,0.308,0.5659999999999999,0.591,0.722,0.309)
    15 = @(7.011,7.935,8.643000000000001,7.803,8.253,8.105,8.93,7.144,7.256,8.507999999999999,7.715,5.153)
}

# Column order matching the $newValues arrays: B..M (2..13)
$cols = @(2,3,4,5,6,7,8,9,10,11,12,13)

foreach ($r in $newValues.Keys) {
    $rowVals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}

Write-Host "Updated sector breakdown trades/closing prices for rows 2-15 (B:M)"
